$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.399145483970642
$ws.Range("B1").Value = 2.47265100479126
$ws.Range("C1").Value = 4.365824699401855
$ws.Range("D1").Value = 4.44134521484375
$ws.Range("E1").Value = 1.486837863922119
